# Regenerate the "K" column (column G) values in the save_data sheet.
# This mirrors the upstream change: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals" — the K column values
# were recomputed and rewritten for the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 2
    11 = 1
    13 = 0
    14 = 2
    15 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
